$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 33 (shifts existing rows 33:99 down to 34:100).
$ws.Rows("33:33").Insert()

# Populate the newly inserted row 33 with the new weekly record
# (Hortaliza, Macroferia Regional de Talca - Sandia).
$ws.Cells.Item(33, 1).Value = 5
$ws.Cells.Item(33, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(33, 3).Value = "Maule"
$ws.Cells.Item(33, 4).Value = 44533
$ws.Cells.Item(33, 5).Value = 7
$ws.Cells.Item(33, 6).Value = 100112028
$ws.Cells.Item(33, 7).Value = "Sandia"
$ws.Cells.Item(33, 8).Value = "Sin especificar"
$ws.Cells.Item(33, 9).Value = "Primera"
$ws.Cells.Item(33, 10).Value = 2000
$ws.Cells.Item(33, 11).Value = 500
$ws.Cells.Item(33, 12).Value = 500
$ws.Cells.Item(33, 13).Value = 500
$ws.Cells.Item(33, 14).Value = "$/kilo"
$ws.Cells.Item(33, 15).Value = "Perú"
$ws.Cells.Item(33, 16).Value = 500
$ws.Cells.Item(33, 17).Value = 1
$ws.Cells.Item(33, 18).Value = "Hortaliza"
